# Apply the recorded workbook edits:
#  1. Rename the single worksheet from "工作表1" to "D9701212"
#     (Excel automatically keeps the _xlnm._FilterDatabase / _xlnm.Database
#     defined names in sync with the new sheet name.)
#  2. Clear the (redundant, fill-less) fill formatting that was applied to
#     cell B2, returning it to the workbook's default cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename worksheet -> defined names referencing it are updated automatically.
$ws.Name = "D9701212"

# 2) Remove the no-op fill formatting from B2 so it reverts to the default style.
$b2 = $ws.Range("B2")
$b2.Interior.Pattern = -4142   # xlNone
